# Auto update Excel log
# Appends newly-logged sensor rows to several worksheets (ALERTS, PIR,
# Humidity, Temperature, Proximity) exactly as produced by the SeniorConnect
# logging pipeline on 2026-02-04 / 2026-02-06.

function Set-TextCell($ws, $row, $col, $val) {
    # Force every value to be written as literal text (the source log
    # always stores Date/Timestamp/Hour/Value as plain strings, never as
    # Excel dates/times/numbers).
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Add-LogRows($ws, $startRow, $rows) {
    $r = $startRow
    foreach ($row in $rows) {
        $c = 1
        foreach ($val in $row) {
            Set-TextCell $ws $r $c $val
            $c = $c + 1
        }
        $r = $r + 1
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALERTS sheet: one new MINIMAL alert row
# ---------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")
$alertRows = ,@('2026-02-04','14:57:33','14:00','Bathroom','MINIMAL','MINIMAL ALERT: Bathroom occupied, no motion > 20s.')
Add-LogRows $wsAlerts 2 $alertRows

# ---------------------------------------------------------------------
# PIR sheet: new motion-sensor readings (rows 14-35)
# ---------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @('2026-02-04','14:57:07','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-04','14:57:08','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-04','14:57:13','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-04','14:57:18','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-04','14:57:23','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-04','14:57:28','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-04','14:57:33','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-04','14:57:38','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-04','14:57:43','14:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:38:23','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:38:28','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:38:34','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:38:38','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:38:39','09:00','Bathroom','Motion Detected','Active'),
    @('2026-02-06','09:38:47','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:38:50','09:00','Bathroom','Motion Detected','Active'),
    @('2026-02-06','09:38:58','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:38:58','09:00','Bathroom','Motion Detected','Active'),
    @('2026-02-06','09:39:06','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:39:08','09:00','Bathroom','Motion Detected','Active'),
    @('2026-02-06','09:39:16','09:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','09:39:17','09:00','Bathroom','Motion Detected','Active')
)
Add-LogRows $wsPir 14 $pirRows

# ---------------------------------------------------------------------
# Humidity sheet: new readings (rows 13-19)
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @('2026-02-04','14:57:06','14:00','Bathroom','83.2%','Active'),
    @('2026-02-04','14:57:07','14:00','Bathroom','84.1%','Active'),
    @('2026-02-04','14:57:10','14:00','Bathroom','83.2%','Active'),
    @('2026-02-04','14:57:25','14:00','Bathroom','82.7%','Active'),
    @('2026-02-04','14:57:30','14:00','Bathroom','83.3%','Active'),
    @('2026-02-04','14:57:35','14:00','Bathroom','84.3%','Active'),
    @('2026-02-04','14:57:41','14:00','Bathroom','83.5%','Active')
)
Add-LogRows $wsHumidity 13 $humidityRows

# ---------------------------------------------------------------------
# Temperature sheet: new readings (rows 13-19)
# ---------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @('2026-02-04','14:57:07','14:00','Bathroom','23.1C','Active'),
    @('2026-02-04','14:57:07','14:00','Bathroom','23.1C','Active'),
    @('2026-02-04','14:57:10','14:00','Bathroom','23.1C','Active'),
    @('2026-02-04','14:57:26','14:00','Bathroom','23.1C','Active'),
    @('2026-02-04','14:57:31','14:00','Bathroom','23.1C','Active'),
    @('2026-02-04','14:57:36','14:00','Bathroom','23.1C','Active'),
    @('2026-02-04','14:57:41','14:00','Bathroom','23.1C','Active')
)
Add-LogRows $wsTemperature 13 $temperatureRows

# ---------------------------------------------------------------------
# Proximity sheet: one new door-entry row
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
$proximityRows = ,@('2026-02-04','14:57:09','14:00','Bathroom Door','ENTER','User ENTERED Bathroom')
Add-LogRows $wsProximity 2 $proximityRows
